$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.73076433333334
$ws.Range("H2").Value = 53.19229300000001
$ws.Range("I2").Value = 0.004631884691211661
$ws.Range("J2").Value = 0.00463188469121166
$ws.Range("M2").Value = 24.576554
$ws.Range("N2").Value = 73.729662
$ws.Range("O2").Value = 0.07553767049546639
$ws.Range("P2").Value = 0.07553767049546638
$ws.Range("Q2").Value = 435.7610870994408
$ws.Range("R2").Value = 3921.849783894967
$ws.Range("S2").Value = 0.0003498817795777415
$ws.Range("T2").Value = 0.0003498817795777414
$ws.Range("G3").Value = 17.73076433333334
$ws.Range("H3").Value = 53.19229300000001
$ws.Range("I3").Value = 0.004631884691211661
$ws.Range("J3").Value = 0.00463188469121166
$ws.Range("O3").Value = 0.359764849016532
$ws.Range("P3").Value = 0.359764849016532
$ws.Range("Q3").Value = 2075.408477377116
$ws.Range("R3").Value = 18678.67629639404
$ws.Range("S3").Value = 0.001666389296595749
$ws.Range("T3").Value = 0.001666389296595749
$ws.Range("G4").Value = 17.73076433333334
$ws.Range("H4").Value = 53.19229300000001
$ws.Range("I4").Value = 0.004631884691211661
$ws.Range("J4").Value = 0.00463188469121166
$ws.Range("M4").Value = 55.68784966666667
$ws.Range("N4").Value = 167.063549
$ws.Range("O4").Value = 0.1711603033819035
$ws.Range("P4").Value = 0.1711603033819035
$ws.Range("Q4").Value = 987.3881386697622
$ws.Range("R4").Value = 8886.493248027858
$ws.Range("S4").Value = 0.0007927947889777824
$ws.Range("T4").Value = 0.0007927947889777823
$ws.Range("G5").Value = 17.73076433333334
$ws.Range("H5").Value = 53.19229300000001
$ws.Range("I5").Value = 0.004631884691211661
$ws.Range("J5").Value = 0.00463188469121166
$ws.Range("M5").Value = 128.0392633333333
$ws.Range("N5").Value = 384.11779
$ws.Range("O5").Value = 0.3935371771060981
$ws.Range("P5").Value = 0.3935371771060981
$ws.Range("Q5").Value = 2270.234003576942
$ws.Range("R5").Value = 20432.10603219247
$ws.Range("S5").Value = 0.001822818826060388
$ws.Range("T5").Value = 0.001822818826060388
$ws.Range("I6").Value = 0.9353873458333681
$ws.Range("J6").Value = 0.935387345833368
$ws.Range("M6").Value = 24.576554
$ws.Range("N6").Value = 73.729662
$ws.Range("O6").Value = 0.07553767049546639
$ws.Range("P6").Value = 0.07553767049546638
$ws.Range("Q6").Value = 87999.90367911836
$ws.Range("R6").Value = 791999.1331120653
$ws.Range("S6").Value = 0.07065698111518982
$ws.Range("T6").Value = 0.07065698111518981
$ws.Range("I7").Value = 0.9353873458333681
$ws.Range("J7").Value = 0.935387345833368
$ws.Range("O7").Value = 0.359764849016532
$ws.Range("P7").Value = 0.359764849016532
$ws.Range("Q7").Value = 419118.9886175733
$ws.Range("R7").Value = 3772070.89755816
$ws.Range("S7").Value = 0.3365194872457163
$ws.Range("T7").Value = 0.3365194872457162
$ws.Range("I8").Value = 0.9353873458333681
$ws.Range("J8").Value = 0.935387345833368
$ws.Range("M8").Value = 55.68784966666667
$ws.Range("N8").Value = 167.063549
$ws.Range("O8").Value = 0.1711603033819035
$ws.Range("P8").Value = 0.1711603033819035
$ws.Range("Q8").Value = 199398.3943706628
$ws.Range("R8").Value = 1794585.549335965
$ws.Range("S8").Value = 0.1601011818924328
$ws.Range("T8").Value = 0.1601011818924328
$ws.Range("I9").Value = 0.9353873458333681
$ws.Range("J9").Value = 0.935387345833368
$ws.Range("M9").Value = 128.0392633333333
$ws.Range("N9").Value = 384.11779
$ws.Range("O9").Value = 0.3935371771060981
$ws.Range("P9").Value = 0.3935371771060981
$ws.Range("Q9").Value = 458463.0880504486
$ws.Range("R9").Value = 4126167.792454037
$ws.Range("S9").Value = 0.3681096955800292
$ws.Range("T9").Value = 0.3681096955800292
$ws.Range("G10").Value = 227.2177583333333
$ws.Range("H10").Value = 681.653275
$ws.Range("I10").Value = 0.0593570833501536
$ws.Range("J10").Value = 0.05935708335015359
$ws.Range("M10").Value = 24.576554
$ws.Range("N10").Value = 73.729662
$ws.Range("O10").Value = 0.07553767049546639
$ws.Range("P10").Value = 0.07553767049546638
$ws.Range("Q10").Value = 5584.229507438117
$ws.Range("R10").Value = 50258.06556694306
$ws.Range("S10").Value = 0.004483695803675837
$ws.Range("T10").Value = 0.004483695803675835
$ws.Range("G11").Value = 227.2177583333333
$ws.Range("H11").Value = 681.653275
$ws.Range("I11").Value = 0.0593570833501536
$ws.Range("J11").Value = 0.05935708335015359
$ws.Range("O11").Value = 0.359764849016532
$ws.Range("P11").Value = 0.359764849016532
$ws.Range("Q11").Value = 26596.1271037305
$ws.Range("R11").Value = 239365.1439335745
$ws.Range("S11").Value = 0.02135459212952972
$ws.Range("T11").Value = 0.02135459212952971
$ws.Range("G12").Value = 227.2177583333333
$ws.Range("H12").Value = 681.653275
$ws.Range("I12").Value = 0.0593570833501536
$ws.Range("J12").Value = 0.05935708335015359
$ws.Range("M12").Value = 55.68784966666667
$ws.Range("N12").Value = 167.063549
$ws.Range("O12").Value = 0.1711603033819035
$ws.Range("P12").Value = 0.1711603033819035
$ws.Range("Q12").Value = 12653.26836766367
$ws.Range("R12").Value = 113879.415308973
$ws.Range("S12").Value = 0.01015957639407722
$ws.Range("T12").Value = 0.01015957639407722
$ws.Range("G13").Value = 227.2177583333333
$ws.Range("H13").Value = 681.653275
$ws.Range("I13").Value = 0.0593570833501536
$ws.Range("J13").Value = 0.05935708335015359
$ws.Range("M13").Value = 128.0392633333333
$ws.Range("N13").Value = 384.11779
$ws.Range("O13").Value = 0.3935371771060981
$ws.Range("P13").Value = 0.3935371771060981
$ws.Range("Q13").Value = 29092.79439325136
$ws.Range("R13").Value = 261835.1495392623
$ws.Range("S13").Value = 0.02335921902287082
$ws.Range("T13").Value = 0.02335921902287082
$ws.Range("G14").Value = 2.387458333333333
$ws.Range("H14").Value = 7.162374999999999
$ws.Range("I14").Value = 0.0006236861252666267
$ws.Range("J14").Value = 0.0006236861252666266
$ws.Range("M14").Value = 24.576554
$ws.Range("N14").Value = 73.729662
$ws.Range("O14").Value = 0.07553767049546639
$ws.Range("P14").Value = 0.07553767049546638
$ws.Range("Q14").Value = 58.67549865191666
$ws.Range("R14").Value = 528.07948786725
$ws.Range("S14").Value = [double]"4.711179702298463E-05"
$ws.Range("T14").Value = [double]"4.711179702298461E-05"
$ws.Range("G15").Value = 2.387458333333333
$ws.Range("H15").Value = 7.162374999999999
$ws.Range("I15").Value = 0.0006236861252666267
$ws.Range("J15").Value = 0.0006236861252666266
$ws.Range("O15").Value = 0.359764849016532
$ws.Range("P15").Value = 0.359764849016532
$ws.Range("Q15").Value = 279.4550291929305
$ws.Range("R15").Value = 2515.095262736375
$ws.Range("S15").Value = 0.0002243803446902538
$ws.Range("T15").Value = 0.0002243803446902538
$ws.Range("G16").Value = 2.387458333333333
$ws.Range("H16").Value = 7.162374999999999
$ws.Range("I16").Value = 0.0006236861252666267
$ws.Range("J16").Value = 0.0006236861252666266
$ws.Range("M16").Value = 55.68784966666667
$ws.Range("N16").Value = 167.063549
$ws.Range("O16").Value = 0.1711603033819035
$ws.Range("P16").Value = 0.1711603033819035
$ws.Range("Q16").Value = 132.9524207520972
$ws.Range("R16").Value = 1196.571786768875
$ws.Range("S16").Value = 0.0001067503064157197
$ws.Range("T16").Value = 0.0001067503064157197
$ws.Range("G17").Value = 2.387458333333333
$ws.Range("H17").Value = 7.162374999999999
$ws.Range("I17").Value = 0.0006236861252666267
$ws.Range("J17").Value = 0.0006236861252666266
$ws.Range("M17").Value = 128.0392633333333
$ws.Range("N17").Value = 384.11779
$ws.Range("O17").Value = 0.3935371771060981
$ws.Range("P17").Value = 0.3935371771060981
$ws.Range("Q17").Value = 305.6884062390278
$ws.Range("R17").Value = 2751.19565615125
$ws.Range("S17").Value = 0.0002454436771376686
$ws.Range("T17").Value = 0.0002454436771376685
